$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# Row 5 "Eligible Patients" (or similar) row L5:V5 -- turn the flat 5000
# literals into a 2% annual growth formula chain starting at M5.
$ws.Range("M5").Formula = "=L5*1.02"
$ws.Range("N5:V5").Formula = "=M5*1.02"

# Key assumption inputs on the right-hand assumptions panel.
$ws.Range("Y31").Value = -0.03
$ws.Range("Y32").Value = 0.09

# View-state tweaks captured in the diff (frozen pane anchor + selection).
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Application.ActiveWindow.ScrollColumn = 14
$ws.Range("R24").Select()
